$wb = $excel.ActiveWorkbook

# --- Update existing "data" sheet: refresh the query timestamp in F2 ---
$ws1 = $wb.Worksheets.Item("data")
$ws1.Range("F2").Value = "2021-10-05 14:22:39.128233"

# --- Add the new "metadata" sheet right after "data" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Reuse the existing bold/bordered/centered header style from "data"
# by copying formats (not values) onto the new sheet's header + A2 cells.
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("B1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Severe familial anorexia"
$ws2.Range("C2").Value = 262
# Force "1.3" to stay text (not be coerced to the number 1.3)
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "1.3"
$ws2.Range("E2").Value = "2020-05-07T14:28:41.057208Z"
$ws2.Range("F2").Value = "2021-10-05 14:22:39.124701"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/262/?format=json"

# Keep the first sheet active/selected like the original workbook
$ws1.Select()
